# Implemented and debugged time variable covariates.
# Adds three new time-varying covariate columns (WT, CRCL, DIAL) to the
# NONMEM-style dataset sheet, fills the existing data rows with the
# placeholder "." value used elsewhere in the sheet, and consolidates the
# duplicate H:MM time format that B3 was using onto the shared H:MM style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the added covariate columns.
$ws.Range("G1").Value = "WT"
$ws.Range("H1").Value = "CRCL"
$ws.Range("I1").Value = "DIAL"

# Populate the new columns on the existing data rows with the same "."
# placeholder already used for missing values elsewhere in the sheet.
$ws.Range("G2:I2").Value = "."
$ws.Range("G3:I3").Value = "."

# B3 was pointing at a duplicate H:MM number format (numFmtId 167); re-apply
# the equivalent H:MM format so it collapses onto the shared style (166)
# instead of the redundant one.
$ws.Range("B3").NumberFormat = "H:MM"

# Match the author's final selection state on the sheet.
$ws.Range("G6:H6").Select()
